$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '34.465.30'
$ws.Range('E2').Value = '  +0.31%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.806.37'
$ws.Range('E3').Value = '  +0.30%  '

# Row 4
$ws.Range('E4').Value = '  -0.17%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '224.83'
$ws.Range('E5').Value = '  -1.13%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.595'
$ws.Range('E6').Value = '  +3.46%  '

# Row 7
$ws.Range('E7').Value = '  -0.12%  '

# Row 8
$ws.Range('E8').Value = '  +5.99%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.288'
$ws.Range('E9').Value = '  -3.90%  '

# Row 10
$ws.Range('E10').Value = '  -3.11%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0975'
$ws.Range('E11').Value = '  +1.17%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.066.39'
$ws.Range('E12').Value = '  +0.27%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.08'
$ws.Range('E13').Value = '  -5.21%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.804.22'
$ws.Range('E14').Value = '  +0.39%  '

# Row 15
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.628'
$ws.Range('E15').Value = '  -2.24%  '

# Row 16
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '34.421.09'
$ws.Range('E16').Value = '  +0.21%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.38'
$ws.Range('E17').Value = '  -2.82%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '67.98'
$ws.Range('E18').Value = '  -1.56%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '241.66'
$ws.Range('E19').Value = '  -1.44%  '

# Row 20
$ws.Range('E20').Value = '  -3.23%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.07'
$ws.Range('E21').Value = '  -4.83%  '

# Row 22
$ws.Range('E22').Value = '  -0.16%  '

# Row 23
$ws.Range('E23').Value = '  -1.65%  '

# Row 24
$ws.Range('E24').Value = '  +1.31%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '171.52'
$ws.Range('E25').Value = '  -0.06%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.71'
$ws.Range('E26').Value = '  -3.31%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.38'
$ws.Range('E27').Value = '  +3.12%  '

# Row 28
$ws.Range('E28').Value = '  +1.27%  '

# Row 29
$ws.Range('E29').Value = '  -0.15%  '

# Row 30
$ws.Range('E30').Value = '  -1.50%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.76'
$ws.Range('E31').Value = '  -1.68%  '

# Row 32
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.84'
$ws.Range('E32').Value = '  -4.09%  '

# Row 33
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0513'
$ws.Range('E33').Value = '  -3.44%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.81'
$ws.Range('E34').Value = '  -0.08%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.317.91'
$ws.Range('E35').Value = '  -5.76%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.640'
$ws.Range('E36').Value = '  -4.45%  '

# Row 37
$ws.Range('E37').Value = '  -1.38%  '

# Row 38
$ws.Range('E38').Value = '  -1.08%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.31'
$ws.Range('E39').Value = '  -6.44%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '82.98'
$ws.Range('E40').Value = '  +0.41%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.44'
$ws.Range('E41').Value = '  +0.99%  '

# Row 42
$ws.Range('B42').Value = 'WEMIXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.21'
$ws.Range('E42').Value = '  -1.09%  '

# Row 43
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.81'
$ws.Range('E43').Value = '  -0.58%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.945'
$ws.Range('E44').Value = '  -1.78%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.69'

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0511'
$ws.Range('E46').Value = '  +0.69%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.966.89'
$ws.Range('E47').Value = '  +0.26%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.77'
$ws.Range('E48').Value = '  -4.37%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '102.19'
$ws.Range('E50').Value = '  -1.96%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0₆0119'
$ws.Range('E51').Value = '  -7.64%  '
